{"js": "// The resume's body text lives in a single paragraph at the end of the\n// document: one run containing alternating <w:t> text segments and\n// manual line breaks (<w:br/>) that lay the whole resume out as\n// line-broken plain text. The edit rewrites that paragraph's content\n// end-to-end (new header, reorganized/renamed sections, new project &\n// experience bullets) while keeping the same \"segment, break, segment\u2026\"\n// shape, so we rebuild the paragraph from a list of the new segments.\n\nconst newSegments = [\n  \"ABDUL RAHMAN \",\n  \" Contact Info: \",\n  \" [437-878-5622, arlnu@uwaterloo.ca, www.linkedin.com/in/abdul-rahman-381b852a7] \",\n  \" \",\n  \" EDUCATION \",\n  \" University of Waterloo, Bachelor of Mathematics (Honours) - Anticipated Graduation: 2024 \",\n  \" \",\n  \" SKILLS \",\n  \" Python, SQL, Java Script, Racket \",\n  \" Google Collab, Jupyter Notebook, Tableau, Power BI, Excel, PowerPoint, Word, HTML \",\n  \" Strong communicator, Team Player, Leadership \",\n  \" \",\n  \" CERTIFICATIONS \",\n  \" Introduction to Data Management by Meta \",\n  \" PROJECTS \",\n  \" Data Governance Framework: Python-based framework for data classification, lineage mapping, and quality checks \",\n  \" Data Quality Monitoring and Anomaly Detection Framework: Machine learning-based anomaly detection system \",\n  \" Functions in Python: Pizza ordering service program design \",\n  \" \",\n  \" EXPERIENCE \",\n  \" House Captain \\u2013 Merryland International School, UAE (2022 - 2024) \",\n  \" Led a team of 100+ students, Organized events, Fostered a positive and inclusive team culture. \",\n  \" \",\n  \" Sports Secretary \\u2013 Merryland International School, UAE (2022 - 2024) \",\n  \" Managed and trained school teams, Organized practice sessions, Promoted sportsmanship and teamwork. \",\n  \" \",\n  \" Varsity Tennis Athlete \\u2013 University of Waterloo (2023 -Placed 4th in the 2024 OUA season, Ranked top 5 in the UAE in both U16 and U18.\",\n];\n\nfunction escapeXml(text) {\n  return text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\")\n    .replace(/\"/g, \"&quot;\")\n    .replace(/'/g, \"&apos;\");\n}\n\n// Rebuild the run's contents: <w:t>segment</w:t> pieces joined by <w:br/>,\n// matching how Word represents an in-paragraph line break.\nconst runInnerXml = newSegments\n  .map((segment) => `<w:t xml:space=\"preserve\">${escapeXml(segment)}</w:t>`)\n  .join(\"<w:br/>\");\n\nconst paragraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  `<w:body><w:p><w:r>${runInnerXml}</w:r></w:p></w:body>` +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The resume content paragraph is the very last paragraph in the body.\nconst targetParagraph = paragraphs.items[paragraphs.items.length - 1];\n\ntargetParagraph.insertOoxml(paragraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The resume content paragraph is the last paragraph in the document body\n# (a single run containing many alternating text segments and manual line\n# breaks that render the whole resume as line-broken text). Replace its\n# contents wholesale with the rewritten resume text, joining the new\n# segments with manual line breaks (Chr(11), i.e. <w:br/>) exactly as the\n# original paragraph was structured.\n$targetParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$targetRange = $targetParagraph.Range\n\n$segments = @(\n    'ABDUL RAHMAN ',\n    ' Contact Info: ',\n    ' [437-878-5622, arlnu@uwaterloo.ca, www.linkedin.com/in/abdul-rahman-381b852a7] ',\n    ' ',\n    ' EDUCATION ',\n    ' University of Waterloo, Bachelor of Mathematics (Honours) - Anticipated Graduation: 2024 ',\n    ' ',\n    ' SKILLS ',\n    ' Python, SQL, Java Script, Racket ',\n    ' Google Collab, Jupyter Notebook, Tableau, Power BI, Excel, PowerPoint, Word, HTML ',\n    ' Strong communicator, Team Player, Leadership ',\n    ' ',\n    ' CERTIFICATIONS ',\n    ' Introduction to Data Management by Meta ',\n    ' PROJECTS ',\n    ' Data Governance Framework: Python-based framework for data classification, lineage mapping, and quality checks ',\n    ' Data Quality Monitoring and Anomaly Detection Framework: Machine learning-based anomaly detection system ',\n    ' Functions in Python: Pizza ordering service program design ',\n    ' ',\n    ' EXPERIENCE ',\n    ' House Captain \u2013 Merryland International School, UAE (2022 - 2024) ',\n    ' Led a team of 100+ students, Organized events, Fostered a positive and inclusive team culture. ',\n    ' ',\n    ' Sports Secretary \u2013 Merryland International School, UAE (2022 - 2024) ',\n    ' Managed and trained school teams, Organized practice sessions, Promoted sportsmanship and teamwork. ',\n    ' ',\n    ' Varsity Tennis Athlete \u2013 University of Waterloo (2023 -Placed 4th in the 2024 OUA season, Ranked top 5 in the UAE in both U16 and U18.'\n)\n\n$targetRange.Text = [string]::Join([string][char]11, $segments)\n"}
